$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 62
$ws.Range("F6").Value = 113
$ws.Range("F7").Value = 597
$ws.Range("F8").Value = 313
$ws.Range("F9").Value = 512
$ws.Range("F11").Value = 10583
$ws.Range("F12").Value = 189
$ws.Range("F15").Value = 2016
$ws.Range("F19").Value = 194
$ws.Range("F21").Value = 228
$ws.Range("F22").Value = 1142
$ws.Range("F23").Value = 112
$ws.Range("F24").Value = 178
$ws.Range("F25").Value = 672
$ws.Range("F27").Value = 206
$ws.Range("F28").Value = 2361
$ws.Range("F29").Value = 659
$ws.Range("F30").Value = 3044
$ws.Range("F31").Value = 991
$ws.Range("F32").Value = 716
$ws.Range("F36").Value = 913
$ws.Range("F37").Value = 17
$ws.Range("F38").Value = 19
$ws.Range("F39").Value = 222
$ws.Range("F41").Value = 1178
$ws.Range("F42").Value = 103
$ws.Range("F45").Value = 215
$ws.Range("F46").Value = 47
$ws.Range("F48").Value = 4064

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 4070
$ws.Range("F8").Value = 69
$ws.Range("F14").Value = 269
$ws.Range("F22").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 736
$ws.Range("F3").Value = 403

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 736
$ws.Range("F3").Value = 403
$ws.Range("F8").Value = 62
$ws.Range("F9").Value = 113
$ws.Range("F10").Value = 597
$ws.Range("F11").Value = 313
$ws.Range("F13").Value = 10583
$ws.Range("F14").Value = 189
$ws.Range("F17").Value = 2016
$ws.Range("F21").Value = 194
$ws.Range("F22").Value = 228
$ws.Range("F23").Value = 1142
$ws.Range("F24").Value = 112
$ws.Range("F25").Value = 178
$ws.Range("F26").Value = 4070
$ws.Range("F28").Value = 672
$ws.Range("F30").Value = 206
$ws.Range("F31").Value = 2361
$ws.Range("F32").Value = 659
$ws.Range("F33").Value = 3044
$ws.Range("F34").Value = 991
$ws.Range("F36").Value = 913
$ws.Range("F37").Value = 19
$ws.Range("F39").Value = 222
$ws.Range("F40").Value = 1178
$ws.Range("F42").Value = 103
$ws.Range("F44").Value = 215
$ws.Range("F47").Value = 4064
